$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.1169995834814548;   C = 1.626987699542094;    D = 3.223369029078222;   E = 13.86384647080068;  G = 18.83120278290246 }
    3 = @{ B = 3.272327238179451;    C = 1.626987699542094;    D = 0.7210945179870265;  E = 0.5333859586016987; G = 6.15379541431027 }
    4 = @{ B = 0.01253208636536152;  C = 0.00006708468553440206; D = 0.1496068669990043; E = 0.5333859586016987; G = 0.6955919966515989 }
    5 = @{ B = 3.272327238179451;    C = 1.626987699542094;    D = 0.1496068669990043;  E = 0.5333859586016987; G = 5.582307763322248 }
    6 = @{ B = 0.1169995834814548;   C = 0.3048912486333797;   D = 0.1496068669990043;  E = 0.5333859586016987; G = 1.104883657715537 }
    7 = @{ B = 3.272327238179451;    C = 1.626987699542094;    D = 0.1496068669990043;  E = 0.5333859586016987; G = 5.582307763322248 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
